# Remove trailing periods from the bullet-point sentences in column E
# (rows 2-14) of the teaching experience worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 5)   # column E
    $val = $cell.Value2
    if ($val -ne $null -and $val.ToString().EndsWith(".")) {
        $cell.Value2 = $val.ToString().Substring(0, $val.ToString().Length - 1)
    }
}

# Move the active selection, matching the saved cursor position in the diff.
$ws.Range("E23").Select()
